$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver License")
$ws.Range("C25").Value = "Required (63-1-35(1)) NOTE: It's named driver license #"
$chars = $ws.Range("C25").Characters(26, 4)
$chars.Font.Bold = $true
Write-Host "done"
